$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "85.799.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.247.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.88%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.360"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +25.87%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.646"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.246.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.571"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.56%  "
$ws.Range("E12").Value = "  +7.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.864.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "85.734.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.260.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "425.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.432.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "75.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000127"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.173"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +18.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "542.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.81%  "
$ws.Range("E38").Value = "  -10.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.388"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "157.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "175.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.96%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.732"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.67%  "
